# Refresh market-data derived figures (prices/profits) on the Leve profit sheets.
# Mirrors a scheduled market-board data pull: only numeric H:N columns change,
# per-row, across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61: Not Taking No for an Answer (Mega-Potion of Strength)
$ws.Range("H61").Value = 4106096.8
$ws.Range("I61").Value = 5714368.5
$ws.Range("J61").Value = 85417
$ws.Range("K61").Value = 17143105.5
$ws.Range("L61").Value = 256251
$ws.Range("M61").Value = -17142933.5
$ws.Range("N61").Value = -256595

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 1268.0769
$ws.Range("I61").Value = 767.8333
$ws.Range("J61").Value = 1696.8572
$ws.Range("K61").Value = 767.8333
$ws.Range("L61").Value = 1696.8572
$ws.Range("M61").Value = -555.8333
$ws.Range("N61").Value = -2120.8572

# Row 101: Art Imitates Life (Doman Steel Tabard of Fending)
$ws.Range("H101").Value = 29999.8
$ws.Range("J101").Value = 29999.8
$ws.Range("L101").Value = 29999.8
$ws.Range("N101").Value = -36489.8

# Row 110: Scheduled Maintenance (Deepgold Ingot)
$ws.Range("H110").Value = 43483380
$ws.Range("I110").Value = 47624556
$ws.Range("K110").Value = 47624556
$ws.Range("M110").Value = -47622511

# Row 135: Forgiveness for My Shins (Ruthenium Sabatons of Fending)
$ws.Range("H135").Value = 28082.25
$ws.Range("J135").Value = 28082.25
$ws.Range("L135").Value = 28082.25
$ws.Range("N135").Value = -38222.25

# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 1268.0769
$ws.Range("I136").Value = 767.8333
$ws.Range("J136").Value = 1696.8572
$ws.Range("K136").Value = 2303.4999
$ws.Range("L136").Value = 5090.571599999999
$ws.Range("M136").Value = 246.5001000000002
$ws.Range("N136").Value = -10190.5716

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 1970.3846
$ws.Range("I134").Value = 2060.147
$ws.Range("J134").Value = 1360
$ws.Range("K134").Value = 6180.441
$ws.Range("L134").Value = 4080
$ws.Range("M134").Value = -3645.441
$ws.Range("N134").Value = -9150

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 21431.766
$ws.Range("I31").Value = 26869.41
$ws.Range("J31").Value = 3759.4167
$ws.Range("K31").Value = 26869.41
$ws.Range("L31").Value = 3759.4167
$ws.Range("M31").Value = -26574.41
$ws.Range("N31").Value = -4349.4167

# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 21431.766
$ws.Range("I34").Value = 26869.41
$ws.Range("J34").Value = 3759.4167
$ws.Range("K34").Value = 26869.41
$ws.Range("L34").Value = 3759.4167
$ws.Range("M34").Value = -26667.41
$ws.Range("N34").Value = -4163.4167

# Row 105: Zelkova, My Love (Zelkova Lumber)
$ws.Range("H105").Value = 1066.238
$ws.Range("I105").Value = 1005
$ws.Range("J105").Value = 1147.8889
$ws.Range("K105").Value = 1005
$ws.Range("L105").Value = 1147.8889
$ws.Range("M105").Value = 742
$ws.Range("N105").Value = -4641.8889

# Row 141: No Greater Treasure (Claro Walnut Necklace of Gathering)
$ws.Range("H141").Value = 73533.336
$ws.Range("J141").Value = 35300
$ws.Range("L141").Value = 35300
$ws.Range("N141").Value = -45660

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap (Maple Syrup)
$ws.Range("H5").Value = 15828.071
$ws.Range("I5").Value = 3142
$ws.Range("J5").Value = 17942.416
$ws.Range("K5").Value = 9426
$ws.Range("L5").Value = 53827.24800000001
$ws.Range("M5").Value = -9314
$ws.Range("N5").Value = -54051.24800000001

# Row 34: Fever Pitch (Chamomile Tea)
$ws.Range("H34").Value = 2633.6155
$ws.Range("J34").Value = 2840.5833
$ws.Range("L34").Value = 8521.749899999999
$ws.Range("N34").Value = -8689.749899999999

# Row 87: Soup That Eats Like a Knight (Clam Chowder)
$ws.Range("H87").Value = 6125
$ws.Range("I87").Value = 6125
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 18375
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -17127
$ws.Range("N87").ClearContents()

# Row 88: Don't Let It Fall Apart (Liver-cheese Sandwich)
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -14572
$ws.Range("N88").ClearContents()

# Row 90: Like Ma Used to Make (L) (Clam Chowder)
$ws.Range("H90").Value = 6125
$ws.Range("I90").Value = 6125
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 55125
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -48885
$ws.Range("N90").ClearContents()

# Row 91: Better Come Back with a Sandwich (L) (Liver-cheese Sandwich)
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -13518
$ws.Range("N91").ClearContents()

# Row 131: The Mountain Steeped (Tsai tou Vounou)
$ws.Range("H131").Value = 794.4
$ws.Range("J131").Value = 833.9
$ws.Range("L131").Value = 2501.7
$ws.Range("N131").Value = -12581.7

# Row 135: Not-so-secret Ingredient (Royal Maple Syrup)
$ws.Range("H135").Value = 15828.071
$ws.Range("I135").Value = 3142
$ws.Range("J135").Value = 17942.416
$ws.Range("K135").Value = 28278
$ws.Range("L135").Value = 161481.744
$ws.Range("M135").Value = -25743
$ws.Range("N135").Value = -166551.744

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit (Mythrite Ingot)
$ws.Range("H70").Value = 60954.223
$ws.Range("I70").Value = 95772.55
$ws.Range("J70").Value = 6239.7144
$ws.Range("K70").Value = 95772.55
$ws.Range("L70").Value = 6239.7144
$ws.Range("M70").Value = -95502.55
$ws.Range("N70").Value = -6779.7144

# Row 73: Hulls of Broken Dreams (L) (Mythrite Ingot)
$ws.Range("H73").Value = 60954.223
$ws.Range("I73").Value = 95772.55
$ws.Range("J73").Value = 6239.7144
$ws.Range("K73").Value = 95772.55
$ws.Range("L73").Value = 6239.7144
$ws.Range("M73").Value = -94836.55
$ws.Range("N73").Value = -8111.7144

# Row 80: Needs More Prayerbell (Hardsilver Ingot)
$ws.Range("H80").Value = 142860460
$ws.Range("J80").Value = 3943.3333
$ws.Range("L80").Value = 3943.3333
$ws.Range("N80").Value = -5939.3333

# Row 83: With a Noise That Reaches Heaven (L) (Hardsilver Ingot)
$ws.Range("H83").Value = 142860460
$ws.Range("J83").Value = 3943.3333
$ws.Range("L83").Value = 19716.6665
$ws.Range("N83").Value = -29700.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck (Dragon Leather)
$ws.Range("H82").Value = 1459.2142
$ws.Range("I82").Value = 1716.5
$ws.Range("J82").Value = 1266.25
$ws.Range("K82").Value = 1716.5
$ws.Range("L82").Value = 1266.25
$ws.Range("M82").Value = -1355.5
$ws.Range("N82").Value = -1988.25

# Row 85: Training Is Only Skintight (L) (Dragon Leather)
$ws.Range("H85").Value = 1459.2142
$ws.Range("I85").Value = 1716.5
$ws.Range("J85").Value = 1266.25
$ws.Range("K85").Value = 1716.5
$ws.Range("L85").Value = 1266.25
$ws.Range("M85").Value = -468.5
$ws.Range("N85").Value = -3762.25

# Row 104: Brace Yourselves (Gazelleskin Bracers of Fending)
$ws.Range("H104").Value = 21200
$ws.Range("J104").Value = 21200
$ws.Range("L104").Value = 21200
$ws.Range("N104").Value = -28188

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 1350.5
$ws.Range("I122").Value = 968
$ws.Range("J122").Value = 1580
$ws.Range("K122").Value = 2904
$ws.Range("L122").Value = 4740
$ws.Range("M122").Value = -454
$ws.Range("N122").Value = -9640
